# Atualizando o arquivo XLSX
# Update odds values in rows 2, 3, 5, 6, 7 of the FlashScore workbook per
# the upstream data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 2.18   # G2: 1.95 -> 2.18
$ws.Cells.Item(2, 8).Value = 2.72   # H2: 2.75 -> 2.72
$ws.Cells.Item(2, 9).Value = 3.95   # I2: 4.85 -> 3.95
$ws.Cells.Item(2, 10).Value = 3   # J2: 2.7 -> 3
$ws.Cells.Item(2, 11).Value = 1.75   # K2: 1.78 -> 1.75
$ws.Cells.Item(2, 12).Value = 4.8   # L2: 5.7 -> 4.8
$ws.Cells.Item(2, 13).Value = 1.18   # M2: 1.17 -> 1.18
$ws.Cells.Item(2, 14).Value = 4.3   # N2: 4.4 -> 4.3
$ws.Cells.Item(2, 15).Value = 1.78   # O2: 1.72 -> 1.78
$ws.Cells.Item(2, 16).Value = 1.93   # P2: 2 -> 1.93
$ws.Cells.Item(2, 17).Value = 3.25   # Q2: 3.1 -> 3.25
$ws.Cells.Item(2, 18).Value = 1.3   # R2: 1.32 -> 1.3
$ws.Cells.Item(2, 19).Value = 6.1   # S2: 5.8 -> 6.1
$ws.Cells.Item(2, 20).Value = 1.09   # T2: 1.1 -> 1.09
$ws.Cells.Item(2, 22).Value = 2   # V2: 2.02 -> 2
$ws.Cells.Item(2, 23).Value = 2.6   # W2: 2.55 -> 2.6
$ws.Cells.Item(2, 24).Value = 1.44   # X2: 1.45 -> 1.44
$ws.Cells.Item(2, 25).Value = 4.6   # Y2: 4.45 -> 4.6
$ws.Cells.Item(2, 26).Value = 8.25   # Z2: 7.4 -> 8.25
$ws.Cells.Item(2, 27).Value = 10.75   # AA2: 10 -> 10.75
$ws.Cells.Item(2, 28).Value = 21   # AB2: 17.5 -> 21
$ws.Cells.Item(2, 29).Value = 28   # AC2: 23 -> 28
$ws.Cells.Item(2, 30).Value = 70   # AD2: 60 -> 70
$ws.Cells.Item(2, 31).Value = 4.3   # AE2: 4.4 -> 4.3
$ws.Cells.Item(2, 35).Value = 6.9   # AI2: 8 -> 6.9
$ws.Cells.Item(2, 36).Value = 18.5   # AJ2: 25 -> 18.5
$ws.Cells.Item(2, 37).Value = 16   # AK2: 19 -> 16
$ws.Cells.Item(2, 38).Value = 70   # AL2: 110 -> 70
$ws.Cells.Item(2, 39).Value = 60   # AM2: 90 -> 60
$ws.Cells.Item(2, 40).Value = 100   # AN2: 120 -> 100
$ws.Cells.Item(3, 7).Value = 2.02   # G3: 1.95 -> 2.02
$ws.Cells.Item(3, 8).Value = 3.7   # H3: 3.75 -> 3.7
$ws.Cells.Item(3, 9).Value = 3.1   # I3: 3.2 -> 3.1
$ws.Cells.Item(3, 10).Value = 2.52   # J3: 2.5 -> 2.52
$ws.Cells.Item(3, 12).Value = 3.5   # L3: 3.6 -> 3.5
$ws.Cells.Item(3, 15).Value = 1.21   # O3: 1.2 -> 1.21
$ws.Cells.Item(3, 16).Value = 3.55   # P3: 3.6 -> 3.55
$ws.Cells.Item(3, 19).Value = 2.5   # S3: 2.47 -> 2.5
$ws.Cells.Item(3, 20).Value = 1.4   # T3: 1.42 -> 1.4
$ws.Cells.Item(3, 26).Value = 10.75   # Z3: 10.5 -> 10.75
$ws.Cells.Item(3, 27).Value = 8.75   # AA3: 8.5 -> 8.75
$ws.Cells.Item(3, 28).Value = 18.5   # AB3: 17.5 -> 18.5
$ws.Cells.Item(3, 29).Value = 15   # AC3: 14.5 -> 15
$ws.Cells.Item(3, 30).Value = 23   # AD3: 22 -> 23
$ws.Cells.Item(3, 31).Value = 13   # AE3: 13.5 -> 13
$ws.Cells.Item(3, 32).Value = 7.3   # AF3: 7.4 -> 7.3
$ws.Cells.Item(3, 35).Value = 11.75   # AI3: 12 -> 11.75
$ws.Cells.Item(3, 36).Value = 17.5   # AJ3: 18.5 -> 17.5
$ws.Cells.Item(3, 37).Value = 11   # AK3: 11.5 -> 11
$ws.Cells.Item(3, 41).Value = 350   # AO3: 300 -> 350
$ws.Cells.Item(5, 7).Value = 1.3   # G5: 1.33 -> 1.3
$ws.Cells.Item(5, 10).Value = 1.69   # J5: 1.73 -> 1.69
$ws.Cells.Item(5, 11).Value = 2.62   # K5: 2.63 -> 2.62
$ws.Cells.Item(6, 7).Value = 1.24   # G6: 1.27 -> 1.24
$ws.Cells.Item(6, 10).Value = 1.63   # J6: 1.67 -> 1.63
$ws.Cells.Item(6, 11).Value = 2.62   # K6: 2.63 -> 2.62
$ws.Cells.Item(6, 14).Value = 11   # N6: 12 -> 11
$ws.Cells.Item(6, 15).Value = 1.14   # O6: 1.13 -> 1.14
$ws.Cells.Item(6, 16).Value = 5   # P6: 5.5 -> 5
$ws.Cells.Item(7, 7).Value = 1.95   # G7: 1.91 -> 1.95
$ws.Cells.Item(7, 8).Value = 3.45   # H7: 3.5 -> 3.45
$ws.Cells.Item(7, 9).Value = 3.6   # I7: 3.7 -> 3.6
$ws.Cells.Item(7, 10).Value = 2.5   # J7: 2.47 -> 2.5
$ws.Cells.Item(7, 12).Value = 4.05   # L7: 4.1 -> 4.05
$ws.Cells.Item(7, 14).Value = 6.9   # N7: 7 -> 6.9
$ws.Cells.Item(7, 19).Value = 3.4   # S7: 3.35 -> 3.4
$ws.Cells.Item(7, 20).Value = 1.27   # T7: 1.28 -> 1.27
$ws.Cells.Item(7, 21).Value = 1.39   # U7: 1.38 -> 1.39
$ws.Cells.Item(7, 22).Value = 2.77   # V7: 2.8 -> 2.77
$ws.Cells.Item(7, 25).Value = 6.8   # Y7: 6.7 -> 6.8
$ws.Cells.Item(7, 26).Value = 8.75   # Z7: 8.5 -> 8.75
$ws.Cells.Item(7, 28).Value = 16.5   # AB7: 16 -> 16.5
$ws.Cells.Item(7, 31).Value = 6.9   # AE7: 7 -> 6.9
$ws.Cells.Item(7, 32).Value = 6.6   # AF7: 6.7 -> 6.6
$ws.Cells.Item(7, 35).Value = 9.75   # AI7: 10 -> 9.75
$ws.Cells.Item(7, 36).Value = 18.5   # AJ7: 19 -> 18.5
